$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    # Two-step find + assign avoids the Find/Replace engine's
    # smart-quote autocorrect substitution that Execute()'s built-in
    # replacement performs on straight apostrophes/quotes.
    $r = $d.Content
    $r.Find.Execute($old) | Out-Null
    $r.Text = $new
}

# ---------------------------------------------------------------
# Title
# ---------------------------------------------------------------
Replace-Text "The Elusive Enigma of Quantum Consciousness" "Unveiling the Enigmatic Symphony of Life: An Exploration of Biological Diversity"

# ---------------------------------------------------------------
# Author name
# ---------------------------------------------------------------
Replace-Text " Alexandro Clemente" " Sarah Jones"

# ---------------------------------------------------------------
# Email (net visible text becomes sarah.jones@valid.com)
# ---------------------------------------------------------------
Replace-Text "aclemente94@protonmail" "sarah.jones@valid"

# ---------------------------------------------------------------
# Main body paragraph
# ---------------------------------------------------------------
Replace-Text "Within the vast tapestry of scientific inquiry, there lies a poignant enigma that has captivated the minds of researchers and philosophers alike: the intricate relationship between quantum mechanics and human consciousness" "Before the dawn of humanity, life existed in a diverse primordial soup, and as time swept over the earth, countless organisms emerged, each a unique melody in the grand symphony of life"

Replace-Text " Quantum mechanics, probing the realm of the exceedingly small, has revealed a puzzling world governed by probability and uncertainty, challenging our conventional understanding of reality" " Biological diversity, the intricate tapestry of life forms, underpins the health of our planet and unveils enigmatic mysteries about the symphony we inhabit"

# Collapse several sentences (+ the two intervening <w:br/>) into one new sentence.
$rngDelStart = $d.Content
$rngDelStart.Find.Execute(" Simultaneously, the enigmatic nature of consciousness, the subjective experience of the world, eludes our grasp") | Out-Null
$startDel = $rngDelStart.Start
$rngDelEnd = $d.Content
$rngDelEnd.Find.Execute(" Concurrently, philosophers and neuroscientists have grappled with the subjective nature of consciousness, exploring how the intricate workings of the brain give rise to our rich experiences of the world.") | Out-Null
$endDel = $rngDelEnd.End
$rngWhole = $d.Range($startDel, $endDel)
$rngWhole.Text = " From the tiniest microbes to the majestic whales, the exploration of biological diversity enriches our understanding of ecology, evolution, and the interdependence of all living organisms"

Replace-Text "As we embark on this exploration, we encounter a profound paradox" "Our planet hosts an array of ecosystems, each a finely tuned ensemble of organisms, shaped by climatic, geographic, and biotic factors"

Replace-Text " The laws of quantum mechanics appear to operate in a realm far removed from our everyday perceptions, seemingly devoid of conscious observers" " Within these ecosystems, diversity plays a pivotal role"

Replace-Text " Yet, it is through our conscious awareness that we attempt to comprehend and interpret the quantum world" " Diverse communities are more resilient to environmental fluctuations, as different species respond to changes in various ways, leading to the adaptation and survival of the ecosystem as a whole"

Replace-Text " This inherent paradox has fueled speculation and research into the potential role of consciousness in shaping or influencing quantum events" " Biological diversity holds the key to understanding the delicate balance of nature, inspiring us to protect and preserve these invaluable ecosystems"

# " The intersection...intriguing " collapses to just a manual line break (no text).
$rngBr = $d.Content
$rngBr.Find.Execute(" The intersection of quantum mechanics and human consciousness thus emerges as an intriguing ") | Out-Null
$rngBr.Text = ""
$rngBr.InsertBreak(6)

# The run that held <w:lastRenderedPageBreak/> + "tapestry of interconnectedness..." becomes
# a plain <w:br/> followed by new text (the rendered-page-break marker disappears), and
# several brand-new sentences are appended after it (before the final, untouched "." run).
$rngTap = $d.Content
$rngTap.Find.Execute("tapestry of interconnectedness, challenging our fundamental assumptions about reality, observation, and the nature of consciousness itself") | Out-Null
$rngTap.Text = "The enigma of biological diversity extends beyond the intricate dance of organisms in ecosystems"
$brPos = $rngTap.Start
$rngTap.InsertAfter(". Evolution, the driving force behind the symphony of life, has led to the remarkable diversity of species and adaptations. By comparing organisms across time and space, scientists piece together the evolutionary puzzle, unraveling the remarkable history of life. Furthermore, the exploration of biological diversity leads to significant scientific advancements with direct implications for human well-being, such as the development of medicines and agricultural innovations")
$rngBreakIns = $d.Range($brPos, $brPos)
$rngBreakIns.InsertBreak(6)

# ---------------------------------------------------------------
# Summary paragraph
# ---------------------------------------------------------------
Replace-Text "The enigmatic connection between quantum mechanics and human consciousness presents a compelling area of research, intertwining the complexities of quantum phenomena with the subjective nature of consciousness" "Biological diversity, the enigmatic symphony of life, unveils the intricate interplay of organisms, unveiling the beauty of ecological balance"

Replace-Text " As we delve into the mysteries of quantum consciousness, we encounter a captivating paradox, compelling us to question the role of conscious observers in the quantum realm" " Exploring this diversity unlocks the secrets of evolution and inspires us to protect our planet's delicate ecosystems"

Replace-Text " Through ongoing explorations, we endeavor to unveil the intricate relationship between these two profound aspects of reality, potentially leading to a deeper understanding of the fundamental nature of matter, consciousness, and the universe we inhabit" " The symphony of life holds invaluable lessons and offers practical benefits that enhance our well-being, making its exploration both fascinating and essential"

# ---------------------------------------------------------------
# Trailing empty paragraph at the end of the document body.
# ---------------------------------------------------------------
$d.Content.InsertParagraphAfter() | Out-Null

Write-Output "done"
